$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric value corrections (recomputed / re-rounded figures) ---
$ws.Range("AD59").Value = -413632.032
$ws.Range("AD60").Value = 625486.0159999999
$ws.Range("AD66").Value = 113165.992
$ws.Range("AH59").Value = -332652.064
$ws.Range("AH63").Value = -165348
$ws.Range("AH64").Value = -124617.992
$ws.Range("AH69").Value = 120343.992
$ws.Range("K61").Value = -224821.968
$ws.Range("K68").Value = 152038
$ws.Range("K73").Value = 95462
$ws.Range("O58").Value = 777193.856
$ws.Range("O59").Value = -295961.152
$ws.Range("O64").Value = -156670.992
$ws.Range("V63").Value = -265701.008
$ws.Range("V66").Value = -131853
$ws.Range("V67").Value = 47574
$ws.Range("V68").Value = -184078
$ws.Range("Z58").Value = 683779.968
$ws.Range("Z73").Value = 126655.992
$ws.Range("Z79").Value = 13905

# --- Clear column R (31/12/2002) for rows 58-79: data removed/blanked ---
$ws.Range("R58:R79").ClearContents()

# --- Clear row 78 (Part. de Acionistas Nao Controladores) except label A78 and already-empty C78 ---
$ws.Range("B78").ClearContents()
$ws.Range("D78:AI78").ClearContents()
